$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, "Bitcoin", "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc", "68.519.48", "  -1.55%  "),
    @(3, "Ethereum", "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth", "2.450.10", "  -2.28%  "),
    @(4, "TetherUSD", "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt", "1.00", "  +0.02%  "),
    @(5, "BNB", "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb", "563.25", "  -2.05%  "),
    @(6, "Solana", "https://coinranking.com/coin/zNZHO_Sjf+solana-sol", "162.89", "  -2.15%  "),
    @(7, "USDC", "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc", "1.00", "  +0.03%  "),
    @(8, "XRP", "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp", "0.505", "  -1.64%  "),
    @(9, "LidoStakedEther", "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth", "2.449.02", "  -2.28%  "),
    @(10, "Dogecoin", "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge", "0.152", "  -5.91%  "),
    @(11, "TRON", "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx", "0.164", "  -1.89%  "),
    @(12, "Cardano", "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada", "0.342", "  -3.81%  "),
    @(13, "Toncoin", "https://coinranking.com/coin/67YlI0K1b+toncoin-ton", "4.82", "  -1.97%  "),
    @(14, "WrappedliquidstakedEther2.0", "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth", "2.903.02", "  -2.10%  "),
    @(15, "WrappedBTC", "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc", "68.440.39", "  -1.52%  "),
    @(16, "ShibaInu", "https://coinranking.com/coin/xz24e0BjL+shibainu-shib", "0.0000171", "  -3.66%  "),
    @(17, "Avalanche", "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax", "23.71", "  -4.67%  "),
    @(18, "WrappedEther", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth", "2.436.98", "  -3.32%  "),
    @(19, "Chainlink", "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link", "10.90", "  -3.07%  "),
    @(20, "BitcoinCash", "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch", "348.63", "  -0.30%  "),
    @(21, "Uniswap", "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni", "7.19", "  -4.52%  "),
    @(22, "Polkadot", "https://coinranking.com/coin/25W7FG7om+polkadot-dot", "3.81", "  -2.51%  "),
    @(23, "SuiNetwork", "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui", "1.88", "  -3.45%  "),
    @(24, "Dai", "https://coinranking.com/coin/MoTuySvg7+dai-dai", "1.00", "  +0.01%  "),
    @(25, "Litecoin", "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc", "68.29", "  -2.61%  "),
    @(26, "NEARProtocol", "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near", "3.77", "  -4.72%  "),
    @(27, "WrappedeETH", "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth", "2.582.90", "  -2.18%  "),
    @(28, "Binance-PegBSC-USD", "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd", "1.01", "  +1.70%  "),
    @(29, "Aptos", "https://coinranking.com/coin/HGYj5JCv5+aptos-apt", "8.28", "  -6.09%  "),
    @(30, "PEPE", "https://coinranking.com/coin/03WI8NQPF+pepe-pepe", "0.0₃0841", "  -5.46%  "),
    @(31, "InternetComputer(DFINITY)", "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp", "7.32", "  -6.59%  "),
    @(32, "Fetch.AI", "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet", "1.18", "  -3.08%  "),
    @(33, "FirstDigitalUSD", "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd", "1.00", "  +0.02%  "),
    @(34, "Bittensor", "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao", "430.76", "  -6.52%  "),
    @(35, "PancakeSwap", "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake", "1.68", "  -2.92%  "),
    @(36, "POPCAT", "https://coinranking.com/coin/sLBuDEsp6+popcat-popcat", "2.99", "  +101.69%  "),
    @(37, "Monero", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr", "157.37", "  -1.58%  "),
    @(38, "WhiteBITCoin", "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt", "19.00", "  -0.29%  "),
    @(39, "USDe", "https://coinranking.com/coin/exbfr2U-0+usde-usde", "1.00", "  +0.05%  "),
    @(40, "Kaspa", "https://coinranking.com/coin/V8GxkwWow+kaspa-kas", "0.110", "  -5.52%  "),
    @(41, "EthereumClassic", "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc", "17.94", "  -2.90%  "),
    @(42, "PolygonEcosystemToken", "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol", "0.307", "  -3.70%  "),
    @(43, "RenderToken", "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render", "4.51", "  -3.54%  "),
    @(44, "Stacks", "https://coinranking.com/coin/mMPrMcB7+stacks-stx", "1.53", "  -3.96%  "),
    @(45, "ImmutableX", "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx", "1.09", "  +0.31%  "),
    @(46, "dogwifhat", "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif", "2.09", "  -5.77%  "),
    @(47, "Aave", "https://coinranking.com/coin/ixgUfzmLR+aave-aave", "135.36", "  -4.83%  "),
    @(48, "Filecoin", "https://coinranking.com/coin/ymQub4fuB+filecoin-fil", "3.37", "  -2.77%  "),
    @(49, "ARBITRUM", "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb", "0.490", "  -5.37%  "),
    @(50, "Cronos", "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro", "0.0718", "  -2.21%  "),
    @(51, "Mantle", "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt", "0.564", "  -2.81%  ")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]

    $dCell = $ws.Cells.Item($r, 4)
    $dCell.NumberFormat = "@"
    $dCell.Value = $row[3]
    $dCell.Style = "Normal"

    $eCell = $ws.Cells.Item($r, 5)
    $eCell.NumberFormat = "@"
    $eCell.Value = $row[4]
    $eCell.Style = "Normal"
}
